$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D10").Value = -8.099
$ws.Range("D12").Value = -6.725
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("D37").Value = -8.315
$ws.Range("D55").Value = -8.218
$ws.Range("D68").Value = -7.228999999999999
$ws.Range("D77").Value = -7.840000000000001
$ws.Range("D78").Value = -8.279
